$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2399.8333
$ws.Range("I80").Value = 2399.8
$ws.Range("K80").Value = 7199.400000000001
$ws.Range("M80").Value = -6201.400000000001
$ws.Range("H83").Value = 2399.8333
$ws.Range("I83").Value = 2399.8
$ws.Range("K83").Value = 21598.2
$ws.Range("M83").Value = -16606.2
$ws.Range("H98").Value = 329.52173
$ws.Range("I98").Value = 329.52173
$ws.Range("K98").Value = 329.52173
$ws.Range("M98").Value = 1168.47827
$ws.Range("H107").Value = 1401.6
$ws.Range("I107").Value = 1356.7693
$ws.Range("J107").Value = 1484.8572
$ws.Range("K107").Value = 1356.7693
$ws.Range("L107").Value = 1484.8572
$ws.Range("M107").Value = 563.2307000000001
$ws.Range("N107").Value = -5324.8572
$ws.Range("H113").Value = 4118.6665
$ws.Range("I113").Value = 2947.2727
$ws.Range("J113").Value = 5959.4287
$ws.Range("K113").Value = 2947.2727
$ws.Range("L113").Value = 5959.4287
$ws.Range("M113").Value = 306.7273
$ws.Range("N113").Value = -12467.4287
$ws.Range("H122").Value = 329.52173
$ws.Range("I122").Value = 329.52173
$ws.Range("K122").Value = 988.56519
$ws.Range("M122").Value = 1461.43481
$ws.Range("H125").Value = 6621
$ws.Range("I125").Value = 10750
$ws.Range("J125").Value = 5441.2856
$ws.Range("K125").Value = 96750
$ws.Range("L125").Value = 48971.5704
$ws.Range("M125").Value = -94290
$ws.Range("N125").Value = -53891.5704
$ws.Range("H137").Value = 12156.581
$ws.Range("I137").Value = 14594.72
$ws.Range("K137").Value = 43784.16
$ws.Range("M137").Value = -41234.16
$ws.Range("H138").Value = 16670905
$ws.Range("I138").Value = 903.53845
$ws.Range("J138").Value = 29418552
$ws.Range("K138").Value = 2710.61535
$ws.Range("L138").Value = 88255656
$ws.Range("M138").Value = 2429.38465
$ws.Range("N138").Value = -88265936

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 30000
$ws.Range("I23").Value = 30000
$ws.Range("K23").Value = 30000
$ws.Range("M23").Value = -29741
$ws.Range("H74").Value = 40658.965
$ws.Range("I74").Value = 46564.348
$ws.Range("J74").Value = 2274
$ws.Range("K74").Value = 46564.348
$ws.Range("L74").Value = 2274
$ws.Range("M74").Value = -45690.348
$ws.Range("N74").Value = -4022
$ws.Range("H77").Value = 40658.965
$ws.Range("I77").Value = 46564.348
$ws.Range("J77").Value = 2274
$ws.Range("K77").Value = 232821.74
$ws.Range("L77").Value = 11370
$ws.Range("M77").Value = -228453.74
$ws.Range("N77").Value = -20106
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 1520.2106
$ws.Range("I122").Value = 1326.8889
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3980.6667
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1530.6667
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 78152.83
$ws.Range("I132").Value = 10339.031
$ws.Range("J132").Value = 801500
$ws.Range("K132").Value = 31017.093
$ws.Range("L132").Value = 2404500
$ws.Range("M132").Value = -28487.093
$ws.Range("N132").Value = -2409560

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3002.7144
$ws.Range("I99").Value = 2652.4546
$ws.Range("K99").Value = 2652.4546
$ws.Range("M99").Value = -1154.4546

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7038.8667
$ws.Range("I22").Value = 7038.8667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 7038.8667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -6688.8667
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 3214.394
$ws.Range("I31").Value = 1388.6666
$ws.Range("J31").Value = 6409.4165
$ws.Range("K31").Value = 1388.6666
$ws.Range("L31").Value = 6409.4165
$ws.Range("M31").Value = -1093.6666
$ws.Range("N31").Value = -6999.4165
$ws.Range("H34").Value = 3214.394
$ws.Range("I34").Value = 1388.6666
$ws.Range("J34").Value = 6409.4165
$ws.Range("K34").Value = 1388.6666
$ws.Range("L34").Value = 6409.4165
$ws.Range("M34").Value = -1186.6666
$ws.Range("N34").Value = -6813.4165
$ws.Range("H62").Value = 3585.5
$ws.Range("I62").Value = 3524.9285
$ws.Range("J62").Value = 3670.3
$ws.Range("K62").Value = 3524.9285
$ws.Range("L62").Value = 3670.3
$ws.Range("M62").Value = -2900.9285
$ws.Range("N62").Value = -4918.3
$ws.Range("H65").Value = 3585.5
$ws.Range("I65").Value = 3524.9285
$ws.Range("J65").Value = 3670.3
$ws.Range("K65").Value = 17624.6425
$ws.Range("L65").Value = 18351.5
$ws.Range("M65").Value = -14504.6425
$ws.Range("N65").Value = -24591.5
$ws.Range("H122").Value = 1409
$ws.Range("I122").Value = 1667.1111
$ws.Range("J122").Value = 247.5
$ws.Range("K122").Value = 5001.3333
$ws.Range("L122").Value = 742.5
$ws.Range("M122").Value = -2551.3333
$ws.Range("N122").Value = -5642.5
$ws.Range("H134").Value = 1812.6666
$ws.Range("I134").Value = 1477.4546
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 4432.3638
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -1897.3638
$ws.Range("N134").Value = -21570
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43940496
$ws.Range("I4").Value = 63637960
$ws.Range("J4").Value = 23011938
$ws.Range("K4").Value = 190913880
$ws.Range("L4").Value = 69035814
$ws.Range("M4").Value = -190913768
$ws.Range("N4").Value = -69036038
$ws.Range("H12").Value = 805
$ws.Range("I12").Value = 322
$ws.Range("K12").Value = 966
$ws.Range("M12").Value = -793
$ws.Range("H23").Value = 624
$ws.Range("I23").Value = 218.16667
$ws.Range("K23").Value = 654.50001
$ws.Range("M23").Value = -419.50001
$ws.Range("H56").Value = 32138.285
$ws.Range("I56").Value = 32138.285
$ws.Range("K56").Value = 32138.285
$ws.Range("M56").Value = -31608.285
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H97").Value = 1165.25
$ws.Range("I97").Value = 1165.25
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3495.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2999.75
$ws.Range("N97").ClearContents()
$ws.Range("H113").Value = 1278
$ws.Range("I113").Value = 427.23077
$ws.Range("J113").Value = 4964.6665
$ws.Range("K113").Value = 1281.69231
$ws.Range("L113").Value = 14893.9995
$ws.Range("M113").Value = 888.3076900000001
$ws.Range("N113").Value = -19233.9995
$ws.Range("H121").Value = 633.5
$ws.Range("I121").Value = 348.57144
$ws.Range("J121").Value = 855.1111
$ws.Range("K121").Value = 1045.71432
$ws.Range("L121").Value = 2565.3333
$ws.Range("M121").Value = 264.28568
$ws.Range("N121").Value = -5185.3333
$ws.Range("H132").Value = 4151.5
$ws.Range("I132").Value = 1389.4
$ws.Range("J132").Value = 5686
$ws.Range("K132").Value = 12504.6
$ws.Range("L132").Value = 51174
$ws.Range("M132").Value = -9974.6
$ws.Range("N132").Value = -56234
$ws.Range("H138").Value = 3806.3
$ws.Range("I138").Value = 3999.625
$ws.Range("K138").Value = 11998.875
$ws.Range("M138").Value = -6858.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4245.2856
$ws.Range("I113").Value = 3487.5
$ws.Range("K113").Value = 3487.5
$ws.Range("M113").Value = -1317.5
$ws.Range("H126").Value = 4299.125
$ws.Range("I126").Value = 3602.75
$ws.Range("K126").Value = 10808.25
$ws.Range("M126").Value = -8338.25
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3801.8823
$ws.Range("I122").Value = 3509.0715
$ws.Range("K122").Value = 10527.2145
$ws.Range("M122").Value = -8077.2145
$ws.Range("H132").Value = 2346.1538
$ws.Range("J132").Value = 4124
$ws.Range("L132").Value = 12372
$ws.Range("N132").Value = -17432

Write-Host "Updated all sheets"